$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Groups")
$ws.Activate()

# Insert a new column at the front, shifting the existing report columns
# (and their widths/styles) one column to the right.
$ws.Columns("A:A").Insert()

# New first column is an "Index" column; give it the same header style as
# the other header cells in row 2.
$ws.Range("B2").Copy()
$ws.Range("A2").PasteSpecial(-4122)
$ws.Range("A2").Value = "Index"

$ws.Application.CutCopyMode = $false

$ws.Range("C13").Select()
